# Actualización automática 2025-08-29 09:20:09
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO": row 19 (INODOROS / LAVABOS columns H,I) ---
$wsVentasGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentasGrupo.Range("H19").Value = 660.6
$wsVentasGrupo.Range("I19").Value = 203.4

# --- Sheet "VENTA MENSUAL": row 19 (F) and totals row 34 (F) ---
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsVentaMensual.Range("F19").Value = 5244.49
$wsVentaMensual.Range("F34").Value = 34255.92

# --- Sheet "CUMPLIMIENTO MENSUAL": rows 7, 8, and totals row 19 ---
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D7").Value = 660.6
$wsCumplimiento.Range("E7").Value = 939.4
$wsCumplimiento.Range("F7").Value = 0.412875

$wsCumplimiento.Range("D8").Value = 203.4
$wsCumplimiento.Range("E8").Value = 421.6
$wsCumplimiento.Range("F8").Value = 0.32544

$wsCumplimiento.Range("D19").Value = 34774.13
$wsCumplimiento.Range("E19").Value = -2664.848924442127
$wsCumplimiento.Range("F19").Value = 1.082993104646951
